# Rename "strain" sheets to "strain_log2_expression" per commit message,
# and switch the active/selected tab from optimization_parameters to
# the dcin5_log2_expression sheet (mirrors the source diff).

$wb = $excel.ActiveWorkbook

$wtSheet = $wb.Worksheets.Item("wt")
$wtSheet.Name = "wt_log2_expression"

$dcin5Sheet = $wb.Worksheets.Item("dcin5")
$dcin5Sheet.Name = "dcin5_log2_expression"

# Make dcin5_log2_expression the active/selected sheet (was optimization_parameters).
$dcin5Sheet.Activate()
